# Update the "All Procedures for Study X" template:
#  - collapse the 3-column header down to a single "Procedure Number" column
#    (drop "Procedure Name" / "Associated CPT Code")
#  - grow the procedure-number list from 25 rows to 100 rows
#  - leave the selection resting on the header row (A1:C1)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 2): keep "Procedure Number" in A2, blank out B2/C2 ---
# Pull B2/C2's formatting from a data row (B3/C3) so they end up with the
# plain bordered style used throughout the data rows, then clear the text.
$ws.Range("B3:C3").Copy() | Out-Null
$ws.Range("B2:C2").PasteSpecial(-4122) | Out-Null
$ws.Range("B2:C2").ClearContents() | Out-Null

# --- Extend the data rows from 25 to 100 procedures (rows 3-102) ---
# Row 27 (procedure #25) is the last existing data row; clone its
# formatting down through row 102 (procedure #100), then fill in the
# sequential procedure numbers for the newly-added rows.
$ws.Range("A27:C27").Copy() | Out-Null
$ws.Range("A28:C102").PasteSpecial(-4122) | Out-Null

for ($i = 28; $i -le 102; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

$excel.CutCopyMode = 0

# --- Selection: leave it on the header row ---
$ws.Range("A1:C1").Select() | Out-Null
